$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header columns (row 1) to the new machine-friendly names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2) Title-case the Spanish connector words ("de", "del", "la", "el", "los",
#    "las", "y") inside the state (col A) and municipality (col B) text, and
#    normalize "TOTAL" to "Total", for every data row.
$connectors = @("de", "del", "la", "el", "los", "las", "y")
for ($r = 2; $r -le 1222; $r++) {
    foreach ($c in 1, 2) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value2
        if ($v -ne $null) {
            if ($v -eq "TOTAL") {
                $cell.Value = "Total"
            } else {
                $words = $v.Split(" ")
                $changed = $false
                for ($i = 0; $i -lt $words.Length; $i++) {
                    if ($connectors -contains $words[$i]) {
                        $words[$i] = $words[$i].Substring(0, 1).ToUpper() + $words[$i].Substring(1)
                        $changed = $true
                    }
                }
                if ($changed) {
                    $cell.Value = [string]::Join(" ", $words)
                }
            }
        }
    }
}

# 3) Remove the trailing footnote/metadata rows (1224-1228) that sat below
#    the blank row 1223; this also shrinks the sheet dimension down to
#    A1:D1222 automatically.
$ws.Range("A1224:D1228").ClearContents() | Out-Null
